# DCTBFSSeeding.xlsx update: new seeding, new .exe build
# - Swap the title/subtitle rows (A1 <-> A2)
# - Add a new "v0.2.3" label at H2
# - Rename the left table's "DEP." header (C3) to "DELTA"
# - Add a second, re-seeded ranking table in columns H:K (rows 3-36)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. Swap A1 (title) and A2 (version) ----
$a1 = $ws.Range("A1").Value2
$a2 = $ws.Range("A2").Value2
$ws.Range("A1").Value = $a2
$ws.Range("A2").Value = $a1

# ---- 2. New version label next to the new table ----
$ws.Range("H2").Value = "v0.2.3"

# ---- 3. Rename C3 header from DEP. to DELTA ----
$ws.Range("C3").Value = "DELTA"

# ---- 4. Build the new seed-delta table in H:K ----
# Header row (copy formatting from the existing A3:D3 header first)
$ws.Range("A3:D3").Copy()
$ws.Range("H3:K3").PasteSpecial(-4122)
$ws.Range("H3").Value = "SEED"
$ws.Range("I3").Value = "NAME"
$ws.Range("J3").Value = "DELTA"
$ws.Range("K3").Value = "V/U"

# Data rows 4-35: copy formatting from corresponding A:D row of the same style group
$ws.Range("A4:D4").Copy()
$ws.Range("H4:K35").PasteSpecial(-4122)

# Row 36 has no A:F counterpart, so copy formatting from row 35 instead
$ws.Range("A35:D35").Copy()
$ws.Range("H36:K36").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# New seeding order (rank, name, delta-from-previous-seed)
$seeding = @(
    @(1,  "Elementus",    1),
    @(2,  "Solcore",     -1),
    @(3,  "Anomal",     $null),
    @(4,  "Martis",      $null),
    @(5,  "Duncan",       3),
    @(6,  "The Demon",   -1),
    @(7,  "Humanus",      2),
    @(8,  "The Rock",    -2),
    @(9,  "Armaments",  $null),
    @(10, "Blizerd",     -1),
    @(11, "Duskan",      -3),
    @(12, "Student",    $null),
    @(13, "Cicle",       $null),
    @(14, "David",       $null),
    @(15, "Doc",         -8),
    @(16, "Spirtu",      $null),
    @(17, "Joseph",       2),
    @(18, "Livern",      -3),
    @(19, "Lyfebud",     $null),
    @(20, "Magnaur",     $null),
    @(21, "Quake",       $null),
    @(22, "Splarg",      $null),
    @(23, "Tweedle",      8),
    @(24, "Zodium",      $null),
    @(25, "Albatross",   -8),
    @(26, "Discrinius",  -1),
    @(27, "Firia",       -1),
    @(28, "Forseer",    -10),
    @(29, "Irode",       -2),
    @(30, "Pyrocitus",   -2),
    @(31, "Rig",         -2),
    @(32, "Stone Golem", -2),
    @(33, "Zip",         -1)
)

$row = 4
foreach ($entry in $seeding) {
    $ws.Cells.Item($row, 8).Value = $entry[0]   # H - SEED
    $ws.Cells.Item($row, 9).Value = $entry[1]   # I - NAME
    if ($entry[2] -ne $null) {
        $ws.Cells.Item($row, 10).Value = $entry[2]   # J - DELTA
    }
    $ws.Cells.Item($row, 11).Value = "V"        # K - V/U
    $row = $row + 1
}

# ---- 5. Column widths for the new I/J columns ----
$ws.Columns.Item(9).ColumnWidth = 17.666666666666664
$ws.Columns.Item(10).ColumnWidth = 15.333333333333332

# ---- 6. View state: scroll / selection ----
$ws.Range("J22").Select()
